$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.122632145881653
$ws.Range("B1").Value = 1.8541339635849
$ws.Range("C1").Value = 6.380153179168701
$ws.Range("D1").Value = 3.29100775718689
$ws.Range("E1").Value = 1.331285953521729
